$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo'd city name in the test data (singa -> singap[ore])
$ws.Range("C2").Value = "singap"

# Move the active selection (the file was re-saved with the view scrolled
# back to the top-left and the cursor resting on C10)
$ws.Range("C10").Select()
